# Generate Report for Handoff
# Regenerate the "Latest Handoff Datetime" (column H) for the
# 4b647b34-46ab-454f-8905-9f77375c347d file row, in both the
# "zh-cn" and "de-de" locale sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H4").Value = "2016-10-19 23:23:44"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H4").Value = "2016-10-19 23:23:54"
